{"js": "// \"update hearing type string\" \u2014 the standalone paragraph that reads\n// \"by <<hearingType>>\" loses its leading \"by \" so it just reads\n// \"<<hearingType>>\" (the bookmark around the placeholder is preserved).\nconst body = context.document.body;\n\nconst hits = body.search(\"by <<hearingType>>\", { matchCase: true });\nhits.load(\"items/text\");\nawait context.sync();\n\nif (hits.items.length > 0) {\n  const hit = hits.items[0];\n  // Narrow the match down to just the \"by \" prefix (3 characters) so the\n  // bookmark + placeholder runs that follow are left alone.\n  const prefixHits = hit.search(\"by \", { matchCase: true });\n  prefixHits.load(\"items/text\");\n  await context.sync();\n\n  if (prefixHits.items.length > 0) {\n    prefixHits.items[0].delete();\n    await context.sync();\n  }\n}\n", "ps1": "# \"update hearing type string\" \u2014 the standalone paragraph that reads\n# \"by <<hearingType>>\" loses its leading \"by \" so it just reads\n# \"<<hearingType>>\" (the bookmark around the placeholder is preserved).\n$d = $word.ActiveDocument\n\n$find = $d.Content.Find\n$find.ClearFormatting()\n$find.MatchWildcards = $false\n$find.Text = \"by <<hearingType>>\"\n$found = $find.Execute()\n\nif ($found) {\n    $hit = $find.Parent\n\n    # Narrow the match down to just the \"by \" prefix (3 characters) so the\n    # bookmark + placeholder runs that follow are left alone.\n    $prefix = $hit.Duplicate\n    $prefix.Start = $hit.Start\n    $prefix.End = $hit.Start + 3\n    if ($prefix.Text -eq \"by \") {\n        $prefix.Delete()\n    }\n}\n"}
